# New questions and lists
#
# In the "OSMI RAZRED" (8th grade) table, the row whose "OBLAST/Subtheme"
# cell reads "Poliedar" has an empty task-number cell (last column).
# We need to fill it in with "173, 174" + ", " + "175" as three separate
# runs, matching the existing formatting used throughout that column
# (w:rFonts w:cstheme="minorHAnsi").

$d = $word.ActiveDocument

# Locate the table / row / column that holds the "Poliedar" subtheme and
# an empty task-number cell right after it.
$targetTable = $null
$targetRow = 0
$targetCol = 0

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $subtheme = $tbl.Cell($r, 3).Range.Text
        $subtheme = $subtheme.TrimEnd([char]7, [char]13)
        if ($subtheme -eq "Poliedar") {
            $numCell = $tbl.Cell($r, 5).Range.Text
            $numCell = $numCell.TrimEnd([char]7, [char]13)
            if ($numCell -eq "") {
                $targetTable = $i
                $targetRow = $r
                $targetCol = 5
            }
        }
    }
}

if ($targetTable -eq $null) {
    throw "Could not find the empty 'Poliedar' task-number cell"
}

# Step 1: put the plain text into the empty cell (collapses the range so
# the paragraph mark / w:p element stays exactly as it was).
$tbl = $d.Tables.Item($targetTable)
$cell = $tbl.Cell($targetRow, $targetCol)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "173, 174, 175"

# Step 2: re-fetch the cell (the previous Range handle is stale after the
# structural edit above) and replace the freshly-typed run with three
# separate runs carrying the correct rFonts formatting, split exactly as
# "173, 174" / ", " / "175" - matching how the surrounding cells in this
# column were authored.
$tbl2 = $d.Tables.Item($targetTable)
$cell2 = $tbl2.Cell($targetRow, $targetCol)
$cellRng2 = $cell2.Range
$runRng = $d.Range($cellRng2.Start, $cellRng2.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="07F8EB77" w14:textId="77777777" w:rsidR="00EF611B" w:rsidRPr="00E14A3A" w:rsidRDefault="00EF611B" w:rsidP="00EF611B">
            <w:pPr>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
              </w:rPr>
              <w:t>173, 174</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
              </w:rPr>
              <w:t xml:space="preserve">, </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:cstheme="minorHAnsi"/>
              </w:rPr>
              <w:t>175</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$runRng.InsertXML($xml)
